$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = 45962
$ws.Range("A24").NumberFormat = $ws.Range("A23").NumberFormat

$ws.Range("B24").Value = 6326
$ws.Range("C24").Value = 1001
$ws.Range("D24").Value = 5923299
$ws.Range("E24").Value = 936.3419222257351
$ws.Range("F24").Value = 7.841800204568705
$ws.Range("G24").Value = 3.730569948186524
$ws.Range("H24").Value = 25.47841659866643
